# Casos de Uso Usuário.docx - apply commit "descrição de casos de uso
# atualizada e protótipos iniciados. close #150"

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $false, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# --- UC01: tidy "Será enviado um email..." sentence (removes spell-check
# markers around "email", no visible text change) ---
Replace-Text "Será enviado um email para o administrador para que ele possa entrar no sistema e atualizar os demais dados do evento e do seu perfil." "Será enviado um email para o administrador para que ele possa entrar no sistema e atualizar os demais dados do evento e do seu perfil."

# --- UC02 header: "UC02" + " - " -> "UC02 - " ---
Replace-Text "UC02 - Manter Gestor Evento" "UC02 - Manter Gestor Evento"

# --- UC02 description: tidy email sentence + password email sentence ---
Replace-Text "dos gestores do evento" "dos gestores do evento"
Replace-Text " Pode-se também enviar email para mudança de senha do gestor no sistema." " Pode-se também enviar email para mudança de senha do gestor no sistema."

# --- UC03 header: "UC03" + " - " -> "UC03 - " ---
Replace-Text "UC03 - Manter Áreas de Interesse" "UC03 - Manter Áreas de Interesse"

# --- UC03 description: tidy run splits ---
Replace-Text "CRUD de áreas de interesse. As áreas " "CRUD de áreas de interesse. As áreas "

# --- UC04 header: "UC04" + " - " -> "UC04 - " ---
Replace-Text "UC04 - Manter Gestor do Sistema" "UC04 - Manter Gestor do Sistema"

# --- UC05 description list items ---
Replace-Text "Cadastrar Participantes" "Manter Colaboradores"
Replace-Text "Cadastrar Colaboradores" "Gerenciar Eventos"
Replace-Text "Manter Eventos" "Gerenciar SubEventos"
Replace-Text "Manter SubEventos" "Adicionar Usuários"

# --- UC05 Protótipos: "UC05 – T01 E UC05 – T02" -> "UC05 – T01, UC05 – T02"
#     with a _GoBack bookmark inserted between the two comma-separated parts ---
Replace-Text "UC05 – T01 E UC05 – T02" "UC05 – T01, UC05 – T02"

Write-Output "done"
